$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 45224 (2023-10-25)
# to 45233 (2023-11-03), matching the automatic update of files.
foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -eq 45224) {
        $cell.Value = 45233
    }
}
